# Auto-generated edit script: updates crypto price/volume table
# to reflect the refreshed GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.337.08'
$ws.Range('E2').Value = '  +3.36%  '
$ws.Range('D3').Value = '1.718.17'
$ws.Range('E3').Value = '  +3.37%  '
$ws.Range('D4').Value = '''0.9995'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '''237.99'
$ws.Range('E5').Value = '  +0.99%  '
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('D7').Value = '''0.4733'
$ws.Range('E7').Value = '  -0.84%  '
$ws.Range('D8').Value = '''0.2624'
$ws.Range('E8').Value = '  +0.61%  '
$ws.Range('D9').Value = '''0.06206'
$ws.Range('E9').Value = '  +0.92%  '
$ws.Range('D10').Value = '1.716.37'
$ws.Range('E10').Value = '  +4.13%  '
$ws.Range('D11').Value = '''0.07051'
$ws.Range('E11').Value = '  -0.24%  '
$ws.Range('D12').Value = '''15.29'
$ws.Range('E12').Value = '  +4.13%  '
$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D13').Value = '''0.5904'
$ws.Range('E13').Value = '  +1.32%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = '''4.414'
$ws.Range('E14').Value = '  +1.69%  '
$ws.Range('D15').Value = '''75.92'
$ws.Range('E16').Value = '  -0.01%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '26.331.91'
$ws.Range('E17').Value = '  +3.41%  '
$ws.Range('B18').Value = 'BinanceUSD'
$ws.Range('C18').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D18').Value = '''1.000'
$ws.Range('E18').Value = '  +0.02%  '
$ws.Range('D19').Value = '''0.000006830'
$ws.Range('E19').Value = '  +1.17%  '
$ws.Range('D20').Value = '''11.53'
$ws.Range('E20').Value = '  +1.63%  '
$ws.Range('D21').Value = '1.934.91'
$ws.Range('E21').Value = '  +2.49%  '
$ws.Range('D22').Value = '''4.536'
$ws.Range('E22').Value = '  +3.08%  '
$ws.Range('D23').Value = '''8.725'
$ws.Range('E23').Value = '  +1.78%  '
$ws.Range('D24').Value = '''5.300'
$ws.Range('E24').Value = '  +0.72%  '
$ws.Range('D25').Value = '''134.76'
$ws.Range('E25').Value = '  +1.41%  '
$ws.Range('D26').Value = '''15.18'
$ws.Range('E26').Value = '  +1.07%  '
$ws.Range('D27').Value = '''1.397'
$ws.Range('E27').Value = '  +0.34%  '
$ws.Range('D28').Value = '''107.35'
$ws.Range('E28').Value = '  +2.89%  '
$ws.Range('D29').Value = '''1.754'
$ws.Range('E29').Value = '  +4.69%  '
$ws.Range('D30').Value = '''3.949'
$ws.Range('E30').Value = '  -0.10%  '
$ws.Range('D31').Value = '''3.679'
$ws.Range('E31').Value = '  +1.19%  '
$ws.Range('D32').Value = '''0.07727'
$ws.Range('E32').Value = '  +1.39%  '
$ws.Range('D33').Value = '''0.04431'
$ws.Range('E33').Value = '  +4.43%  '
$ws.Range('D34').Value = '''2.612'
$ws.Range('D35').Value = '''0.9726'
$ws.Range('E35').Value = '  +3.22%  '
$ws.Range('D36').Value = '''0.6137'
$ws.Range('E36').Value = '  +0.83%  '
$ws.Range('D37').Value = '''0.9242'
$ws.Range('E37').Value = '  +8.52%  '
$ws.Range('D38').Value = '''112.43'
$ws.Range('E38').Value = '  +15.73%  '
$ws.Range('D39').Value = '''2.465'
$ws.Range('E39').Value = '  -5.37%  '
$ws.Range('D40').Value = '''1.919'
$ws.Range('E40').Value = '  +3.69%  '
$ws.Range('D41').Value = '''1.000'
$ws.Range('E41').Value = '  +0.01%  '
$ws.Range('D42').Value = '''0.01471'
$ws.Range('E42').Value = '  -0.75%  '
$ws.Range('D43').Value = '''5.308'
$ws.Range('E43').Value = '  +13.28%  '
$ws.Range('D44').Value = '''0.3811'
$ws.Range('E44').Value = '  +1.82%  '
$ws.Range('D45').Value = '''0.1158'
$ws.Range('E45').Value = '  +3.90%  '
$ws.Range('D46').Value = '''6.271'
$ws.Range('E46').Value = '  +1.68%  '
$ws.Range('D47').Value = '''0.05268'
$ws.Range('E47').Value = '  +0.29%  '
$ws.Range('D48').Value = '''30.20'
$ws.Range('E48').Value = '  +3.09%  '
$ws.Range('D49').Value = '''7.635'
$ws.Range('E49').Value = '  +5.62%  '
$ws.Range('B50').Value = 'Decentraland'
$ws.Range('C50').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D50').Value = '''0.3362'
$ws.Range('E50').Value = '  +1.64%  '
$ws.Range('B51').Value = 'TrueUSD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd'
$ws.Range('D51').Value = '''1.002'
$ws.Range('E51').Value = '  +0.06%  '
